$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B10").Copy()
$ws.Range("B11").PasteSpecial(-4122)

$ws.Range("C10").Copy()
$ws.Range("C11").PasteSpecial(-4122)

$ws.Range("B11").Value = "Requirment 9.1.3 labled as 9.1.2"
$ws.Range("C11").Value = "Fixed"

$excel.CutCopyMode = $false
